# Update "想去人数" (column F) counts that changed between the two data
# pulls, on both the "展览" sheet and the "全部类型" sheet. The two sheets
# list (mostly) the same events but "全部类型" has one extra row inserted
# (row 11: the "东方LiveParty" entry), so from that point on the matching
# rows are offset by 1 between the sheets. Rather than hard-code row
# numbers (which differ per sheet), key the updates by the event name in
# column C so the correct row is updated on each sheet regardless of the
# offset.

$wb = $excel.ActiveWorkbook

# Map of event name (column C) -> new "想去人数" value (column F).
$updates = @{
    '南昌·漫拥出品-晨啼星舟随机宅舞启航场(免费活动)' = 10
    '南昌·SuperComic动漫游戏博览会' = 4587
    '吉安·COMIC LIFE次元假日05' = 703
    '九江·第一届Loading加载中动漫展' = 188
    '赣州·第四届赣州半夏动漫展' = 999
    '南昌·漫拥动漫嘉年华Pro-追光启航' = 239
    '九江·SXD动漫嘉年华' = 71
    '抚州·临次元08·盛夏动漫狂欢节' = 119
    '南昌·萌卡动漫展' = 3544
    '江西·次元星河动漫游戏嘉年华' = 5899
    '赣州·马娘only' = 37
    '万载·第八届馨缘动漫文化展' = 44
    '南昌·幻梦境国际动漫游戏嘉年华1th' = 3367
    '吉安·COMIC LIFE周年庆典' = 364
    '景德镇·第十五届瓷都ACG动漫游戏博览会' = 2475
    '上饶·第十五届IX Group国风嘉年华暨十周年庆典' = 265
    '九江·第一届异次元动漫嘉年华' = 353
    '上饶·囧喵喵国风动漫展' = 128
    '南昌·第一届异次元动漫嘉年华' = 1019
    '赣州·第二届异次元动漫嘉年华' = 913
    '信丰·七夕节UPUP动漫展' = 21
    '南昌·W·MEETING动漫游戏盛典' = 27
    '吉安·WF无线次元新星动漫博览会' = 53
    '上饶·次元重现夏日嘉年华' = 64
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
